# Move "Picture 5" (the networking diagram picture) on slide 2 to the end
# of the shape stack (so it renders on top of / after the other groups)
# and reposition + rotate it into its new spot.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the picture shape by name (it is currently shape #5, right before
# "Group 8" / "Group 7").
$pic = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Picture 5") {
        $pic = $shp
        break
    }
}

# Send it to the very end of the z-order (after "Group 8" and "Group 7").
$pic.ZOrder(0)  # msoBringToFront

# Reposition and rotate the picture into its new location.
$pic.Rotation = 28.07085
$pic.Left = 197.8825196850394
$pic.Top = 312.5994488188977
